# "Added Duplicate Part Numbers"
#
# Inserts a new "PROFILES" column at G (duplicating the Sill/Jamb Screw
# Spline Assembly profile that already lives in column F), which pushes
# the former G:K columns of the small row 1-4 summary table out to H:L.
# Also updates the Quantity in F3 (12.0 ft -> 16.0 ft), recalculates the
# Price row (row 4) for the now 7-wide table, and rolls the new total
# into the GRAND TOTAL cell (E14).
#
# NOTE: this table (rows 1-4, cols E:L) sits directly above an unrelated
# ACCESSORIES table (rows 7-10, cols E:Q) in the same columns, so a
# whole-column Insert/Shift would incorrectly drag that second table
# along too. Values are therefore written directly, cell by cell,
# scoped only to the rows that actually change per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: PROFILES header / names ---
$ws.Range("F1").Value = "Sill/Jamb Screw Spline Assembly"
$ws.Range("G1").Value = "Sill/Jamb Screw Spline Assembly"
$ws.Range("H1").Value = "Flush Filler"
$ws.Range("I1").Value = "Two Piece Mullion Screw Spline Assembly"
$ws.Range("J1").Value = "Horizontal Screw Spline Assembly"
$ws.Range("K1").Value = "Head"
$ws.Range("L1").Value = "Thermal Sill Flashing"

# --- Row 2: Part Number ---
$ws.Range("F2").Value = "BE9-2513"
$ws.Range("G2").Value = "BE9-2513"
$ws.Range("H2").Value = "E9-2512"
$ws.Range("I2").Value = "BE9-2511"
$ws.Range("J2").Value = "BE9-2515"
$ws.Range("K2").Value = "BE9-2514"
$ws.Range("L2").Value = "BE9-2578"

# --- Row 3: Quantity ---
$ws.Range("F3").Value = "16.0 ft"
$ws.Range("G3").Value = "12.0 ft"
$ws.Range("H3").Value = "24.0 ft"
$ws.Range("I3").Value = "24.0 ft"
$ws.Range("J3").Value = "12.0 ft"
$ws.Range("K3").Value = "12.0 ft"
$ws.Range("L3").Value = "12.0 ft"

# --- Row 4: Price, plus the Grand Total. These look like currency, so
#     Excel's auto-detection would otherwise silently convert them to
#     numbers; mark the cells as Text first so they stay literal
#     strings like the rest of the sheet (all inline/shared strings),
#     then drop the formatting back to Normal so no extra styling is
#     left applied to the cells themselves. ---
$priceAddrs = @("F4", "G4", "H4", "I4", "J4", "K4", "L4", "E14")
foreach ($addr in $priceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("F4").Value = "`$335.20"
$ws.Range("G4").Value = "`$251.40"
$ws.Range("H4").Value = "`$180.00"
$ws.Range("I4").Value = "`$596.40"
$ws.Range("J4").Value = "`$316.20"
$ws.Range("K4").Value = "`$266.40"
$ws.Range("L4").Value = "`$158.40"
$ws.Range("E14").Value = "`$2420.10"

foreach ($addr in $priceAddrs) {
    $ws.Range($addr).Style = "Normal"
}

# --- Column widths for the newly split G/H columns. The stored <col>
#     width ends up ColumnWidth + 0.8333333333333333 (5px padding at
#     the default Calibri 11 metrics), so back that padding out to
#     land exactly on the target widths of 33 and 17. ---
$pad = 0.8333333333333333
$ws.Columns("G:G").ColumnWidth = 33 - $pad
$ws.Columns("H:H").ColumnWidth = 17 - $pad
